# Datathon_GRUPO_33.pptx - slide 15 (Dashboard Analitico)
# 1) Reposition/resize the big dashboard screenshot picture ("Imagem 4").
# 2) Add a new centered textbox below it with a hyperlinked Power BI URL.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(15)

# ---- 1) Resize / reposition the picture (id=5, name="Imagem 4") ----
$pic = $s.Shapes.Item(5)
$pic.Left   = 91.16129302978516
$pic.Top    = 136.54397583007812
$pic.Width  = 800.1290283203125
$pic.Height = 341.9076843261719

# ---- 2) Add the new textbox with the Power BI link ----
# A throwaway shape is created first and removed so that the shape-id
# counter lands on 3 for the real textbox (matches the authored file).
$dummy = $s.Shapes.AddTextbox(1, 0, 0, 10, 10)
$dummy.Delete()

$tb = $s.Shapes.AddTextbox(1, 118.45149230957031, 486.5992126464844, 750.193603515625, 41.19842529296875)
$tb.Name = "CaixaDeTexto 2"

$tb.TextFrame.WordWrap = -1
$tb.TextFrame.AutoSize = 1

$tr = $tb.TextFrame.TextRange
$tr.Text = "https://app.powerbi.com/view?r=eyJrIjoiYWE0NGQyNWQtODAzOS00NDIyLWI3MmEtMTU1YzJiZDdhMjNjIiwidCI6IjExZGJiZmUyLTg5YjgtNDU0OS1iZTEwLWNlYzM2NGU1OTU1MSIsImMiOjR9"
$tr.LanguageID = "pt-BR"
$tr.ParagraphFormat.Alignment = 2

$ast = $tr.ActionSettings.Item(1)
$ast.Hyperlink.Address = "https://app.powerbi.com/view?r=eyJrIjoiYWE0NGQyNWQtODAzOS00NDIyLWI3MmEtMTU1YzJiZDdhMjNjIiwidCI6IjExZGJiZmUyLTg5YjgtNDU0OS1iZTEwLWNlYzM2NGU1OTU1MSIsImMiOjR9"

# Restore the exact authored height (turning on AutoSize recalculated it
# once the text was set - Left/Top/Width are untouched by that, only
# Height needs to be pinned back to the committed extent).
$tb.Height = 41.19842529296875

# Keep the shape background transparent (<a:noFill/>).
$tb.Fill.Visible = 0
